$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new header cells (row 1, one per extra table) and the new
# column-A rows (one per extra table's class name), in the exact order
# that reproduces the original author's shared-string append order.
$ws.Range("D1").Value = "TABLE_抽卡1"
$ws.Range("A4").Value = "Table_Gacha1"

$ws.Range("E1").Value = "TABLE_关卡"
$ws.Range("A5").Value = "Table_Quest"

$ws.Range("F1").Value = "TABLE_技能"
$ws.Range("A6").Value = "Table_Skill"

$ws.Range("A7").Value = "Table_Attack"

$ws.Range("G1").Value = "TABLE_英雄攻击"

$ws.Range("H1").Value = "TABLE_英雄成长表"
$ws.Range("A8").Value = "Table_Herogrow"

$ws.Range("I1").Value = "TABLE_英雄经验表"
$ws.Range("A9").Value = "Table_Heroexp"

$ws.Range("J1").Value = "TABLE_角色经验及成长"
$ws.Range("A10").Value = "Table_Role"

$ws.Range("K1").Value = "TABLE_物品合成表"
$ws.Range("A11").Value = "Table_Itemcombine"

$ws.Range("A12").Value = "Table_Enemy1"

$ws.Range("L1").Value = "Enemy/TABLE_Enemy1"

# Resize the newly added columns to fit their (short) header text - mirrors
# the bestFit column widths the workbook ends up with after the new data
# is entered.
$ws.Columns("B").ColumnWidth = 10.857142857142858
$ws.Columns("C").ColumnWidth = 10.857142857142858
$ws.Columns("D").ColumnWidth = 11.714285714285714
$ws.Columns("E").ColumnWidth = 10.857142857142858
$ws.Columns("F").ColumnWidth = 10.857142857142858
$ws.Columns("G").ColumnWidth = 14.857142857142858
$ws.Columns("H").ColumnWidth = 16.857142857142858
$ws.Columns("I").ColumnWidth = 16.857142857142858
$ws.Columns("J").ColumnWidth = 20.857142857142858
$ws.Columns("K").ColumnWidth = 16.857142857142858
$ws.Columns("L").ColumnWidth = 12.857142857142858

$ws.Range("A12").Select() | Out-Null
